# Apply recalculated simulation values (n_samples = 1000 re-run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"39.99962409359842"
$ws.Range("G2").Value = [double]"98.78689035774144"
$ws.Range("H2").Value = [double]"11.58210785781587"
$ws.Range("P2").Value = [double]"38.28134654520498"
$ws.Range("U2").Value = [double]"3.973795235303034e-11"
$ws.Range("V2").Value = [double]"7.208772533205847e-05"
$ws.Range("X2").Value = [double]"0.002651499999998919"
$ws.Range("F3").Value = [double]"39.99986532097269"
$ws.Range("G3").Value = [double]"75.06410640381461"
$ws.Range("H3").Value = [double]"11.64382052896311"
$ws.Range("O3").Value = [double]"19.69541914122515"
$ws.Range("P3").Value = [double]"43.53147100922693"
$ws.Range("U3").Value = [double]"4.996823792972454e-12"
$ws.Range("V3").Value = [double]"2.845094646085731e-05"
$ws.Range("X3").Value = [double]"0.002237299999997333"
$ws.Range("F4").Value = [double]"40.00001491477332"
$ws.Range("G4").Value = [double]"78.47844308854032"
$ws.Range("H4").Value = [double]"3.31113127145463"
$ws.Range("P4").Value = [double]"38.51163498793859"
$ws.Range("U4").Value = [double]"5.619160464056704e-13"
$ws.Range("V4").Value = [double]"3.151370226711346e-05"
$ws.Range("X4").Value = [double]"0.002448599999993917"
$ws.Range("F5").Value = [double]"39.99964695333965"
$ws.Range("G5").Value = [double]"99.14406300181923"
$ws.Range("H5").Value = [double]"11.49762063392832"
$ws.Range("P5").Value = [double]"38.17322568730212"
$ws.Range("U5").Value = [double]"3.529045426769465e-11"
$ws.Range("V5").Value = [double]"6.795416966568627e-05"
$ws.Range("X5").Value = [double]"0.003166399999997793"
$ws.Range("H6").Value = [double]"9.566754714964928"
$ws.Range("P6").Value = [double]"26.07184847509767"
$ws.Range("U6").Value = [double]"2.8037455233838e-13"
$ws.Range("V6").Value = [double]"1.142062456489517e-05"
$ws.Range("X6").Value = [double]"0.001734700000000089"
$ws.Range("F7").Value = [double]"39.99980265404244"
$ws.Range("G7").Value = [double]"39.95891618602371"
$ws.Range("H7").Value = [double]"13.44861595113718"
$ws.Range("P7").Value = [double]"51.07323844922092"
$ws.Range("U7").Value = [double]"1.235321587199202e-11"
$ws.Range("V7").Value = [double]"5.029263453766466e-05"
$ws.Range("X7").Value = [double]"0.001957500000003165"
$ws.Range("H8").Value = [double]"6.98371089718414"
$ws.Range("U8").Value = [double]"1.395623962086032e-12"
$ws.Range("V8").Value = [double]"2.723887120722592e-05"
$ws.Range("X8").Value = [double]"0.001709699999999259"
$ws.Range("H9").Value = [double]"12.35952713437253"
$ws.Range("U9").Value = [double]"3.043174594637796e-11"
$ws.Range("V9").Value = [double]"0.0001319392386701159"
$ws.Range("X9").Value = [double]"0.001460099999995634"
$ws.Range("U10").Value = [double]"5.681223642926786e-11"
$ws.Range("V10").Value = [double]"4.960583444535423e-05"
$ws.Range("X10").Value = [double]"0.002291000000006704"
$ws.Range("F11").Value = [double]"39.99986200749663"
$ws.Range("G11").Value = [double]"1.867343394803042"
$ws.Range("H11").Value = [double]"9.0227607060538"
$ws.Range("P11").Value = [double]"48.22912390133411"
$ws.Range("U11").Value = [double]"2.653443357333249e-12"
$ws.Range("V11").Value = [double]"8.038896756278828e-05"
$ws.Range("X11").Value = [double]"0.001453300000001434"
$ws.Range("G12").Value = [double]"138.8372751052582"
$ws.Range("H12").Value = [double]"3.249002087167057"
$ws.Range("P12").Value = [double]"34.6505162081074"
$ws.Range("U12").Value = [double]"2.798526034180758e-13"
$ws.Range("V12").Value = [double]"2.954063598844842e-05"
$ws.Range("X12").Value = [double]"0.001657599999994375"
$ws.Range("G13").Value = [double]"32.44840812761685"
$ws.Range("H13").Value = [double]"3.233717956305221"
$ws.Range("P13").Value = [double]"40.78814301860535"
$ws.Range("U13").Value = [double]"3.69855338563836e-14"
$ws.Range("V13").Value = [double]"1.581922933295635e-05"
$ws.Range("X13").Value = [double]"0.002305100000000948"
$ws.Range("E14").Value = [double]"7.913927351953512e-06"
$ws.Range("F14").Value = [double]"39.99985964168974"
$ws.Range("G14").Value = [double]"27.74949541810546"
$ws.Range("H14").Value = [double]"11.88039915974985"
$ws.Range("P14").Value = [double]"50.60639292564557"
$ws.Range("U14").Value = [double]"1.108688286374933e-11"
$ws.Range("V14").Value = [double]"6.023459001829533e-05"
$ws.Range("X14").Value = [double]"0.002275900000000775"
$ws.Range("F15").Value = [double]"40.00006897160229"
$ws.Range("G15").Value = [double]"139.9857347344505"
$ws.Range("H15").Value = [double]"5.850455148991594"
$ws.Range("P15").Value = [double]"32.50765448399709"
$ws.Range("U15").Value = [double]"5.608552201974197e-12"
$ws.Range("V15").Value = [double]"4.392012932971592e-05"
$ws.Range("X15").Value = [double]"0.001680200000002685"
$ws.Range("F16").Value = [double]"39.99945625370279"
$ws.Range("H16").Value = [double]"12.27771993193822"
$ws.Range("P16").Value = [double]"40.02678412524576"
$ws.Range("U16").Value = [double]"7.826888616576329e-11"
$ws.Range("V16").Value = [double]"0.0001024049942956816"
$ws.Range("X16").Value = [double]"0.004265300000000138"
$ws.Range("F17").Value = [double]"40.0000070619833"
$ws.Range("H17").Value = [double]"5.146715137565042"
$ws.Range("P17").Value = [double]"32.14042818239013"
$ws.Range("U17").Value = [double]"7.786504539242426e-13"
$ws.Range("V17").Value = [double]"2.878572361584097e-05"
$ws.Range("X17").Value = [double]"0.00160929999999837"
$ws.Range("G18").Value = [double]"40.19311934564482"
$ws.Range("H18").Value = [double]"9.277616524507932"
$ws.Range("P18").Value = [double]"46.64211293629056"
$ws.Range("U18").Value = [double]"5.428821823742404e-11"
$ws.Range("V18").Value = [double]"0.0001413333747860265"
$ws.Range("X18").Value = [double]"0.001564700000002972"
$ws.Range("F19").Value = [double]"39.99991559263758"
$ws.Range("G19").Value = [double]"53.03986458111573"
$ws.Range("O19").Value = [double]"12.95212141825336"
$ws.Range("P19").Value = [double]"45.29432276563003"
$ws.Range("U19").Value = [double]"2.752568494802972e-12"
$ws.Range("V19").Value = [double]"2.779249207878528e-05"
$ws.Range("X19").Value = [double]"0.002268300000004331"
$ws.Range("F20").Value = [double]"39.99997808545956"
$ws.Range("H20").Value = [double]"9.06543493107303"
$ws.Range("P20").Value = [double]"34.2025516523895"
$ws.Range("U20").Value = [double]"1.819168840929539e-13"
$ws.Range("V20").Value = [double]"4.553697306732885e-06"
$ws.Range("X20").Value = [double]"0.001972399999999652"
$ws.Range("F21").Value = [double]"40.00008158930126"
$ws.Range("G21").Value = [double]"17.55591967619605"
$ws.Range("H21").Value = [double]"12.13521265286859"
$ws.Range("P21").Value = [double]"51.5686377853337"
$ws.Range("U21").Value = [double]"8.178432767879544e-12"
$ws.Range("V21").Value = [double]"6.027894357688662e-05"
$ws.Range("X21").Value = [double]"0.001405800000000568"
$ws.Range("P22").Value = [double]"37.03361164013523"
$ws.Range("U22").Value = [double]"1.183538026332665e-12"
$ws.Range("V22").Value = [double]"1.161854545845586e-05"
$ws.Range("X22").Value = [double]"0.003104799999995578"
$ws.Range("F23").Value = [double]"39.99981809915775"
$ws.Range("H23").Value = [double]"6.469799786907129"
$ws.Range("P23").Value = [double]"45.11790376526474"
$ws.Range("U23").Value = [double]"1.562551729961003e-11"
$ws.Range("V23").Value = [double]"0.0002136970288390308"
$ws.Range("X23").Value = [double]"0.00142720000000196"
$ws.Range("F24").Value = [double]"39.99988820127784"
$ws.Range("G24").Value = [double]"37.52320147086804"
$ws.Range("H24").Value = [double]"8.780036733819708"
$ws.Range("P24").Value = [double]"46.36483185253479"
$ws.Range("U24").Value = [double]"7.60329828657267e-12"
$ws.Range("V24").Value = [double]"5.785004208906982e-05"
$ws.Range("X24").Value = [double]"0.001538699999997561"
$ws.Range("F25").Value = [double]"39.99992912262342"
$ws.Range("O25").Value = [double]"26.76285026898852"
$ws.Range("U25").Value = [double]"1.180035594230863e-12"
$ws.Range("V25").Value = [double]"1.209061339148655e-05"
$ws.Range("X25").Value = [double]"0.002270200000005218"
$ws.Range("F26").Value = [double]"39.99999459307674"
$ws.Range("P26").Value = [double]"23.992337447728"
$ws.Range("U26").Value = [double]"1.254084272274159e-12"
$ws.Range("V26").Value = [double]"1.732623755093339e-05"
$ws.Range("X26").Value = [double]"0.001694899999996835"
$ws.Range("H27").Value = [double]"9.987235601957098"
$ws.Range("U27").Value = [double]"5.930608099441129e-11"
$ws.Range("V27").Value = [double]"7.320937541258871e-05"
$ws.Range("X27").Value = [double]"0.001662799999998299"
$ws.Range("G28").Value = [double]"55.42849438202202"
$ws.Range("H28").Value = [double]"9.232624932393245"
$ws.Range("P28").Value = [double]"44.86690715507035"
$ws.Range("U28").Value = [double]"3.908149279947755e-12"
$ws.Range("V28").Value = [double]"3.286978329970637e-05"
$ws.Range("X28").Value = [double]"0.001863499999998908"
$ws.Range("F29").Value = [double]"39.99966381007912"
$ws.Range("H29").Value = [double]"12.77481900422614"
$ws.Range("O29").Value = [double]"21.45486080061464"
$ws.Range("P29").Value = [double]"43.97532707179653"
$ws.Range("U29").Value = [double]"2.885654423952492e-11"
$ws.Range("V29").Value = [double]"6.535815947729094e-05"
$ws.Range("X29").Value = [double]"0.002082899999997778"
$ws.Range("U30").Value = [double]"1.009933979987035e-10"
$ws.Range("X30").Value = [double]"0.001749799999998913"
$ws.Range("F31").Value = [double]"40.0000222356378"
$ws.Range("G31").Value = [double]"155.4764538682488"
$ws.Range("H31").Value = [double]"6.418822550048666"
$ws.Range("O31").Value = [double]"8.182086277360632"
$ws.Range("P31").Value = [double]"30.68830020845654"
$ws.Range("U31").Value = [double]"2.11191379422657e-12"
$ws.Range("V31").Value = [double]"3.357144846777655e-05"
$ws.Range("X31").Value = [double]"0.001565800000001616"
$ws.Range("P32").Value = [double]"34.68993350520547"
$ws.Range("U32").Value = [double]"6.184141198721374e-12"
$ws.Range("V32").Value = [double]"2.818333410833304e-05"
$ws.Range("X32").Value = [double]"0.002052100000000223"
$ws.Range("F33").Value = [double]"39.99993370374865"
$ws.Range("G33").Value = [double]"51.4897251355753"
$ws.Range("H33").Value = [double]"14.49721911908653"
$ws.Range("P33").Value = [double]"50.43112365321306"
$ws.Range("U33").Value = [double]"1.126032976113059e-12"
$ws.Range("V33").Value = [double]"1.358209622632341e-05"
$ws.Range("X33").Value = [double]"0.001963800000005733"
$ws.Range("G34").Value = [double]"154.5995681097097"
$ws.Range("H34").Value = [double]"7.472468171171147"
$ws.Range("P34").Value = [double]"29.7101792478708"
$ws.Range("U34").Value = [double]"7.593578680015849e-12"
$ws.Range("V34").Value = [double]"4.563464655437062e-05"
$ws.Range("X34").Value = [double]"0.001624399999997195"
$ws.Range("G35").Value = [double]"44.73268922387123"
$ws.Range("H35").Value = [double]"8.399400825378075"
$ws.Range("P35").Value = [double]"45.30378643345113"
$ws.Range("U35").Value = [double]"4.944869474893517e-11"
$ws.Range("V35").Value = [double]"0.0001412782789929285"
$ws.Range("X35").Value = [double]"0.001927300000005516"
$ws.Range("H36").Value = [double]"12.36160422767842"
$ws.Range("P36").Value = [double]"52.19292494065064"
$ws.Range("U36").Value = [double]"9.458855497785467e-12"
$ws.Range("V36").Value = [double]"9.646460173848406e-05"
$ws.Range("X36").Value = [double]"0.001464000000005683"
$ws.Range("F37").Value = [double]"39.99995352533437"
$ws.Range("G37").Value = [double]"82.45010274530991"
$ws.Range("H37").Value = [double]"6.700173552578881"
$ws.Range("P37").Value = [double]"39.51189763484432"
$ws.Range("U37").Value = [double]"9.849150890725537e-13"
$ws.Range("V37").Value = [double]"1.73959620969735e-05"
$ws.Range("X37").Value = [double]"0.001855300000002558"
$ws.Range("F38").Value = [double]"39.99997172978126"
$ws.Range("G38").Value = [double]"131.1313524793212"
$ws.Range("H38").Value = [double]"10.21936487128386"
$ws.Range("P38").Value = [double]"31.10941882016252"
$ws.Range("U38").Value = [double]"3.150573095000125e-13"
$ws.Range("V38").Value = [double]"5.027317556991002e-06"
$ws.Range("X38").Value = [double]"0.002004900000002863"
$ws.Range("F39").Value = [double]"39.99999432707852"
$ws.Range("H39").Value = [double]"8.569140498383451"
$ws.Range("U39").Value = [double]"5.063940513815707e-13"
$ws.Range("V39").Value = [double]"1.572664802745939e-05"
$ws.Range("X39").Value = [double]"0.001698100000005809"
$ws.Range("P40").Value = [double]"45.95999675234545"
$ws.Range("U40").Value = [double]"1.759430506423475e-11"
$ws.Range("V40").Value = [double]"5.165794318430166e-05"
$ws.Range("X40").Value = [double]"0.001974900000000446"
$ws.Range("F41").Value = [double]"39.9999780619308"
$ws.Range("P41").Value = [double]"38.67753768299394"
$ws.Range("U41").Value = [double]"1.269215057415164e-13"
$ws.Range("V41").Value = [double]"3.994153780008817e-06"
$ws.Range("X41").Value = [double]"0.002209499999999309"
$ws.Range("G42").Value = [double]"73.38312080831898"
$ws.Range("H42").Value = [double]"3.664412846442649"
$ws.Range("U42").Value = [double]"6.048511016044316e-13"
$ws.Range("V42").Value = [double]"3.069817782830753e-05"
$ws.Range("X42").Value = [double]"0.001691699999994967"
$ws.Range("F43").Value = [double]"40.00000202140297"
$ws.Range("G43").Value = [double]"45.16600088814613"
$ws.Range("H43").Value = [double]"3.553170669750127"
$ws.Range("U43").Value = [double]"1.354798613493489e-13"
$ws.Range("V43").Value = [double]"2.245301442461338e-05"
$ws.Range("X43").Value = [double]"0.001620799999997757"
$ws.Range("F44").Value = [double]"39.99957269140694"
$ws.Range("G44").Value = [double]"67.41747114952439"
$ws.Range("H44").Value = [double]"14.40830771878925"
$ws.Range("P44").Value = [double]"47.23988615183812"
$ws.Range("U44").Value = [double]"4.321251262616544e-11"
$ws.Range("V44").Value = [double]"7.92152248335425e-05"
$ws.Range("X44").Value = [double]"0.001930700000002616"
$ws.Range("G45").Value = [double]"20.70945538460091"
$ws.Range("U45").Value = [double]"2.523254699903326e-11"
$ws.Range("V45").Value = [double]"0.0001193002251490424"
$ws.Range("X45").Value = [double]"0.002116300000004401"
$ws.Range("G46").Value = [double]"150.2712832721844"
$ws.Range("H46").Value = [double]"4.257524893756699"
$ws.Range("U46").Value = [double]"3.821871594658846e-13"
$ws.Range("V46").Value = [double]"2.67204494297785e-05"
$ws.Range("X46").Value = [double]"0.001595000000001789"
$ws.Range("F47").Value = [double]"39.9999883084805"
$ws.Range("P47").Value = [double]"25.92452113144004"
$ws.Range("U47").Value = [double]"8.555195156229673e-14"
$ws.Range("V47").Value = [double]"7.294337981973988e-06"
$ws.Range("X47").Value = [double]"0.001597500000002583"
$ws.Range("E48").Value = [double]"1.481620291598854e-05"
$ws.Range("F48").Value = [double]"39.99967927008353"
$ws.Range("H48").Value = [double]"10.61796520234186"
$ws.Range("U48").Value = [double]"4.597766713631341e-11"
$ws.Range("V48").Value = [double]"5.214465488522199e-05"
$ws.Range("X48").Value = [double]"0.001756100000001481"
$ws.Range("H49").Value = [double]"14.69368008053075"
$ws.Range("U49").Value = [double]"2.683852302159238e-12"
$ws.Range("V49").Value = [double]"2.836979156080521e-05"
$ws.Range("X49").Value = [double]"0.001725199999995652"
$ws.Range("F50").Value = [double]"39.99989925431343"
$ws.Range("H50").Value = [double]"13.69911415780831"
$ws.Range("U50").Value = [double]"2.033794641826449e-12"
$ws.Range("V50").Value = [double]"4.92974291424987e-05"
$ws.Range("X50").Value = [double]"0.001421600000000467"
$ws.Range("G51").Value = [double]"54.17814450998296"
$ws.Range("H51").Value = [double]"5.964054612019525"
$ws.Range("P51").Value = [double]"42.15244758634718"
$ws.Range("U51").Value = [double]"3.426037656060878e-10"
$ws.Range("V51").Value = [double]"0.0004781301554815607"
$ws.Range("X51").Value = [double]"0.001446400000006065"
$ws.Range("G52").Value = [double]"89.03059316754577"
$ws.Range("U52").Value = [double]"6.104686371122455e-11"
$ws.Range("V52").Value = [double]"0.0001469575001053373"
$ws.Range("X52").Value = [double]"0.002355100000002608"
$ws.Range("F53").Value = [double]"40.00001737420603"
$ws.Range("G53").Value = [double]"121.6309094440031"
$ws.Range("H53").Value = [double]"3.536530513231022"
$ws.Range("P53").Value = [double]"35.48196353864436"
$ws.Range("U53").Value = [double]"9.537109387614699e-13"
$ws.Range("V53").Value = [double]"3.369663124452568e-05"
$ws.Range("X53").Value = [double]"0.001679800000005116"
$ws.Range("G54").Value = [double]"148.1595698436984"
$ws.Range("H54").Value = [double]"3.313394664054614"
$ws.Range("P54").Value = [double]"34.16255080386239"
$ws.Range("U54").Value = [double]"1.547789935769679e-13"
$ws.Range("V54").Value = [double]"2.478102446876765e-05"
$ws.Range("X54").Value = [double]"0.001625900000000513"
$ws.Range("F55").Value = [double]"40.00023247943204"
$ws.Range("H55").Value = [double]"9.733939893915139"
$ws.Range("P55").Value = [double]"27.82235863487175"
$ws.Range("U55").Value = [double]"3.885653130005065e-11"
$ws.Range("V55").Value = [double]"5.432692423023478e-05"
$ws.Range("X55").Value = [double]"0.001572899999999322"
$ws.Range("G56").Value = [double]"75.63322759635366"
$ws.Range("H56").Value = [double]"7.35883033207865"
$ws.Range("O56").Value = [double]"13.80603978827702"
$ws.Range("P56").Value = [double]"40.75203632296728"
$ws.Range("U56").Value = [double]"2.845913084737441e-12"
$ws.Range("V56").Value = [double]"2.895087326914095e-05"
$ws.Range("X56").Value = [double]"0.001924400000000048"
$ws.Range("G57").Value = [double]"16.7597770776447"
$ws.Range("H57").Value = [double]"14.65922237913705"
$ws.Range("O57").Value = [double]"6.56737777266306"
$ws.Range("P57").Value = [double]"54.58545535271175"
$ws.Range("U57").Value = [double]"6.090484667048666e-11"
$ws.Range("V57").Value = [double]"0.0001514243460474611"
$ws.Range("X57").Value = [double]"0.001418600000000936"
$ws.Range("F58").Value = [double]"39.9999911533625"
$ws.Range("G58").Value = [double]"177.1545189038948"
$ws.Range("H58").Value = [double]"7.000929995792018"
$ws.Range("P58").Value = [double]"29.03630038828905"
$ws.Range("U58").Value = [double]"3.882866746618497e-14"
$ws.Range("V58").Value = [double]"7.237979845307401e-06"
$ws.Range("X58").Value = [double]"0.001573899999996797"
$ws.Range("G59").Value = [double]"15.72167079643027"
$ws.Range("H59").Value = [double]"5.766369132583212"
$ws.Range("P59").Value = [double]"44.14942243062936"
$ws.Range("U59").Value = [double]"4.110803797188344e-11"
$ws.Range("V59").Value = [double]"0.0003371259004701618"
$ws.Range("X59").Value = [double]"0.001644800000001112"
$ws.Range("F60").Value = [double]"39.99986647708538"
$ws.Range("G60").Value = [double]"65.7193408871492"
$ws.Range("H60").Value = [double]"8.414388990415183"
$ws.Range("P60").Value = [double]"42.81161643432613"
$ws.Range("U60").Value = [double]"6.838659295088161e-12"
$ws.Range("V60").Value = [double]"4.335184597403174e-05"
$ws.Range("X60").Value = [double]"0.001926300000000936"
$ws.Range("H61").Value = [double]"4.894425456490418"
$ws.Range("P61").Value = [double]"31.77996962870389"
$ws.Range("U61").Value = [double]"9.11948647234352e-14"
$ws.Range("V61").Value = [double]"1.427628546116909e-05"
$ws.Range("X61").Value = [double]"0.001673100000004979"
$ws.Range("G62").Value = [double]"54.15431487512916"
$ws.Range("H62").Value = [double]"5.243679128926555"
$ws.Range("P62").Value = [double]"41.54475629238003"
$ws.Range("U62").Value = [double]"6.529499300317126e-11"
$ws.Range("V62").Value = [double]"0.000262070967781524"
$ws.Range("X62").Value = [double]"0.002048399999999617"
$ws.Range("F63").Value = [double]"40.00001744084962"
$ws.Range("H63").Value = [double]"5.176157577505605"
$ws.Range("U63").Value = [double]"1.393447160569871e-12"
$ws.Range("V63").Value = [double]"3.433406352574741e-05"
$ws.Range("X63").Value = [double]"0.001623600000002057"
$ws.Range("G64").Value = [double]"6.075700087764135"
$ws.Range("H64").Value = [double]"8.869691795016239"
$ws.Range("P64").Value = [double]"48.02229189773811"
$ws.Range("U64").Value = [double]"8.872047106201748e-12"
$ws.Range("V64").Value = [double]"0.0001291752313071497"
$ws.Range("X64").Value = [double]"0.001381000000002075"
$ws.Range("F65").Value = [double]"39.9997456844341"
$ws.Range("H65").Value = [double]"7.754890670877845"
$ws.Range("O65").Value = [double]"16.10869069217062"
$ws.Range("U65").Value = [double]"2.560317224463109e-11"
$ws.Range("V65").Value = [double]"6.806732770920982e-05"
$ws.Range("X65").Value = [double]"0.001902300000004686"
$ws.Range("F66").Value = [double]"39.99940928182556"
$ws.Range("G66").Value = [double]"68.23872610698317"
$ws.Range("H66").Value = [double]"9.804836940565798"
$ws.Range("P66").Value = [double]"43.46274857393188"
$ws.Range("U66").Value = [double]"1.136203618370239e-10"
$ws.Range("V66").Value = [double]"0.00015488230219439"
$ws.Range("X66").Value = [double]"0.001957900000000734"
$ws.Range("F67").Value = [double]"39.9999990122092"
$ws.Range("G67").Value = [double]"39.35591441200611"
$ws.Range("H67").Value = [double]"3.376850641112725"
$ws.Range("P67").Value = [double]"40.68916021004919"
$ws.Range("U67").Value = [double]"7.417664636539093e-14"
$ws.Range("V67").Value = [double]"1.931320902539507e-05"
$ws.Range("X67").Value = [double]"0.001602300000001833"
$ws.Range("F68").Value = [double]"39.99999892132351"
$ws.Range("H68").Value = [double]"4.604412405748593"
$ws.Range("P68").Value = [double]"42.4380410789613"
$ws.Range("U68").Value = [double]"5.201257407669949e-14"
$ws.Range("V68").Value = [double]"1.267593946081846e-05"
$ws.Range("X68").Value = [double]"0.00167679999999848"
$ws.Range("F69").Value = [double]"39.99996025664753"
$ws.Range("H69").Value = [double]"13.30626012446609"
$ws.Range("P69").Value = [double]"38.98051049958387"
$ws.Range("U69").Value = [double]"3.985968417207791e-13"
$ws.Range("V69").Value = [double]"6.997352115342007e-06"
$ws.Range("X69").Value = [double]"0.002269100000006574"
$ws.Range("G70").Value = [double]"46.68659596521127"
$ws.Range("H70").Value = [double]"3.520353630061811"
$ws.Range("P70").Value = [double]"40.50452106733947"
$ws.Range("U70").Value = [double]"1.472670033439836e-13"
$ws.Range("V70").Value = [double]"2.32235391329502e-05"
$ws.Range("X70").Value = [double]"0.001555400000000873"
$ws.Range("G71").Value = [double]"53.32487247071196"
$ws.Range("H71").Value = [double]"12.70879281787568"
$ws.Range("U71").Value = [double]"2.265557200516949e-13"
$ws.Range("V71").Value = [double]"6.440196641540065e-06"
$ws.Range("X71").Value = [double]"0.002046200000002329"
$ws.Range("O72").Value = [double]"12.56035885836787"
$ws.Range("U72").Value = [double]"1.314536660778977e-13"
$ws.Range("V72").Value = [double]"2.816427354756055e-06"
$ws.Range("X72").Value = [double]"0.002436699999996961"
$ws.Range("F73").Value = [double]"40.00003571547256"
$ws.Range("G73").Value = [double]"33.40328378613725"
$ws.Range("H73").Value = [double]"7.639284741555365"
$ws.Range("P73").Value = [double]"45.46366636773733"
$ws.Range("U73").Value = [double]"2.418946424356174e-12"
$ws.Range("V73").Value = [double]"4.340871837287285e-05"
$ws.Range("X73").Value = [double]"0.001502100000003281"
$ws.Range("G74").Value = [double]"19.21267023402487"
$ws.Range("H74").Value = [double]"11.21237523053843"
$ws.Range("P74").Value = [double]"50.40105589706384"
$ws.Range("U74").Value = [double]"1.272606611737983e-11"
$ws.Range("V74").Value = [double]"7.943064598815052e-05"
$ws.Range("X74").Value = [double]"0.001393900000003612"
$ws.Range("P75").Value = [double]"34.40678872454651"
$ws.Range("U75").Value = [double]"5.382193603300878e-13"
$ws.Range("V75").Value = [double]"7.870282168344717e-06"
$ws.Range("X75").Value = [double]"0.002039599999996256"
$ws.Range("H76").Value = [double]"12.41452205462801"
$ws.Range("U76").Value = [double]"9.510878389367971e-13"
$ws.Range("V76").Value = [double]"1.522137011674747e-05"
$ws.Range("X76").Value = [double]"0.001917900000002248"
$ws.Range("G77").Value = [double]"102.4084754660749"
$ws.Range("P77").Value = [double]"36.83630882581323"
$ws.Range("U77").Value = [double]"1.285438083454058e-12"
$ws.Range("V77").Value = [double]"2.598318170834181e-05"
$ws.Range("X77").Value = [double]"0.001651999999999987"
$ws.Range("G78").Value = [double]"55.94939953808473"
$ws.Range("H78").Value = [double]"12.82781669040134"
$ws.Range("P78").Value = [double]"48.06395603417501"
$ws.Range("U78").Value = [double]"5.408604161761127e-13"
$ws.Range("V78").Value = [double]"9.762903550879213e-06"
$ws.Range("X78").Value = [double]"0.002035800000001586"
$ws.Range("F79").Value = [double]"40.00008142421624"
$ws.Range("H79").Value = [double]"10.7062538261042"
$ws.Range("U79").Value = [double]"1.515746802521602e-11"
$ws.Range("V79").Value = [double]"4.632511019393969e-05"
$ws.Range("X79").Value = [double]"0.00177569999999605"
$ws.Range("G80").Value = [double]"8.67243420242751"
$ws.Range("P80").Value = [double]"50.73719759273454"
$ws.Range("U80").Value = [double]"1.576381883984814e-11"
$ws.Range("V80").Value = [double]"0.0001272735168826945"
$ws.Range("X80").Value = [double]"0.001449900000004334"
$ws.Range("F81").Value = [double]"39.99939344579809"
$ws.Range("H81").Value = [double]"8.616556174980941"
$ws.Range("P81").Value = [double]"39.09394041447641"
$ws.Range("U81").Value = [double]"1.294161772783958e-10"
$ws.Range("V81").Value = [double]"0.0001557580338348724"
$ws.Range("X81").Value = [double]"0.00186639999999727"
$ws.Range("G82").Value = [double]"93.94518649634435"
$ws.Range("P82").Value = [double]"39.33687975491983"
$ws.Range("U82").Value = [double]"2.337776205767417e-11"
$ws.Range("V82").Value = [double]"5.69317403232753e-05"
$ws.Range("X82").Value = [double]"0.002589000000000397"
$ws.Range("H83").Value = [double]"7.090170627499205"
$ws.Range("P83").Value = [double]"42.05114650405129"
$ws.Range("U83").Value = [double]"2.578697035254993e-13"
$ws.Range("V83").Value = [double]"9.995335842345966e-06"
$ws.Range("X83").Value = [double]"0.00181129999999996"
$ws.Range("H84").Value = [double]"10.67065097312856"
$ws.Range("P84").Value = [double]"46.97776422471183"
$ws.Range("U84").Value = [double]"8.34057889232215e-12"
$ws.Range("V84").Value = [double]"4.486763374909526e-05"
$ws.Range("X84").Value = [double]"0.001804799999995055"
$ws.Range("F85").Value = [double]"39.99986033879235"
$ws.Range("G85").Value = [double]"46.26780636519948"
$ws.Range("H85").Value = [double]"11.2053386203029"
$ws.Range("O85").Value = [double]"13.12798827717038"
$ws.Range("P85").Value = [double]"47.93152832304281"
$ws.Range("U85").Value = [double]"6.795973167088639e-12"
$ws.Range("V85").Value = [double]"4.009325714983593e-05"
$ws.Range("X85").Value = [double]"0.002540299999999718"
$ws.Range("F86").Value = [double]"39.99983111842248"
$ws.Range("H86").Value = [double]"7.44338101375969"
$ws.Range("U86").Value = [double]"1.463457360690266e-11"
$ws.Range("V86").Value = [double]"0.0001827209385627903"
$ws.Range("X86").Value = [double]"0.001522500000000093"
$ws.Range("H87").Value = [double]"9.107428161010303"
$ws.Range("O87").Value = [double]"18.76060995912142"
$ws.Range("P87").Value = [double]"35.2196635194256"
$ws.Range("U87").Value = [double]"3.168057824975849e-13"
$ws.Range("V87").Value = [double]"6.3329697878442e-06"
$ws.Range("X87").Value = [double]"0.002067699999997785"
$ws.Range("G88").Value = [double]"36.4691283801627"
$ws.Range("H88").Value = [double]"5.485141745912376"
$ws.Range("O88").Value = [double]"6.791947957944521"
$ws.Range("P88").Value = [double]"42.97849381140354"
$ws.Range("U88").Value = [double]"1.769625251035028e-13"
$ws.Range("V88").Value = [double]"1.600325847206795e-05"
$ws.Range("X88").Value = [double]"0.001589299999999128"
$ws.Range("F89").Value = [double]"39.99988003120575"
$ws.Range("H89").Value = [double]"14.4587627068685"
$ws.Range("P89").Value = [double]"33.52888714055681"
$ws.Range("U89").Value = [double]"3.742362742370268e-12"
$ws.Range("V89").Value = [double]"1.943379671037936e-05"
$ws.Range("X89").Value = [double]"0.002308199999994542"
$ws.Range("F90").Value = [double]"39.99997685758038"
$ws.Range("H90").Value = [double]"9.085790285801991"
$ws.Range("P90").Value = [double]"34.29462332625654"
$ws.Range("U90").Value = [double]"2.016348103532636e-13"
$ws.Range("V90").Value = [double]"4.817663189591433e-06"
$ws.Range("X90").Value = [double]"0.001955099999996435"
$ws.Range("F91").Value = [double]"39.99939251481509"
$ws.Range("G91").Value = [double]"89.9485858731999"
$ws.Range("H91").Value = [double]"6.071009581688424"
$ws.Range("U91").Value = [double]"1.877962206150732e-10"
$ws.Range("V91").Value = [double]"0.000242414438860229"
$ws.Range("X91").Value = [double]"0.001617000000003088"
$ws.Range("F92").Value = [double]"39.9999995425895"
$ws.Range("G92").Value = [double]"27.54967441079621"
$ws.Range("H92").Value = [double]"4.892246354555277"
$ws.Range("P92").Value = [double]"42.7655545441803"
$ws.Range("U92").Value = [double]"5.829370670662089e-14"
$ws.Range("V92").Value = [double]"1.244867362901912e-05"
$ws.Range("X92").Value = [double]"0.001696099999996648"
$ws.Range("G93").Value = [double]"36.17963127079899"
$ws.Range("H93").Value = [double]"6.948538710436816"
$ws.Range("P93").Value = [double]"44.53291085194562"
$ws.Range("U93").Value = [double]"2.813815998424713e-11"
$ws.Range("V93").Value = [double]"0.0001528363978434941"
$ws.Range("X93").Value = [double]"0.00144809999999751"
$ws.Range("F94").Value = [double]"40.00063064298005"
$ws.Range("P94").Value = [double]"24.24477627832842"
$ws.Range("U94").Value = [double]"2.721804458103254e-10"
$ws.Range("V94").Value = [double]"0.0001062805177442521"
$ws.Range("X94").Value = [double]"0.001616200000000845"
$ws.Range("F95").Value = [double]"39.99978751221068"
$ws.Range("G95").Value = [double]"109.1441024411213"
$ws.Range("P95").Value = [double]"35.75112342661529"
$ws.Range("U95").Value = [double]"1.391470317749134e-11"
$ws.Range("V95").Value = [double]"4.083442914503364e-05"
$ws.Range("X95").Value = [double]"0.002598400000003664"
$ws.Range("F96").Value = [double]"39.99996284139227"
$ws.Range("H96").Value = [double]"10.28664598520448"
$ws.Range("P96").Value = [double]"31.41107720246519"
$ws.Range("U96").Value = [double]"5.297450679295137e-13"
$ws.Range("V96").Value = [double]"6.664405512226042e-06"
$ws.Range("X96").Value = [double]"0.001983199999997964"
$ws.Range("F97").Value = [double]"40.00001033063616"
$ws.Range("G97").Value = [double]"41.18447696115249"
$ws.Range("H97").Value = [double]"5.288281052382939"
$ws.Range("U97").Value = [double]"2.213601620308847e-13"
$ws.Range("V97").Value = [double]"1.746655408260547e-05"
$ws.Range("X97").Value = [double]"0.00161239999999907"
$ws.Range("F98").Value = [double]"39.99980100896806"
$ws.Range("G98").Value = [double]"38.64377338885009"
$ws.Range("H98").Value = [double]"13.91385179683049"
$ws.Range("P98").Value = [double]"51.74065836545003"
$ws.Range("U98").Value = [double]"1.239928576045537e-11"
$ws.Range("V98").Value = [double]"4.988508665312349e-05"
$ws.Range("X98").Value = [double]"0.001809099999995567"
$ws.Range("U99").Value = [double]"1.169119582406338e-11"
$ws.Range("V99").Value = [double]"2.56850620397244e-05"
$ws.Range("X99").Value = [double]"0.002143799999998919"
$ws.Range("P100").Value = [double]"31.16646275320443"
$ws.Range("U100").Value = [double]"1.041550365622056e-10"
$ws.Range("V100").Value = [double]"9.011928441986518e-05"
$ws.Range("X100").Value = [double]"0.001793200000001605"
$ws.Range("H101").Value = [double]"8.856845623041153"
$ws.Range("U101").Value = [double]"8.05782499126468e-11"
$ws.Range("V101").Value = [double]"9.196118749692551e-05"
$ws.Range("X101").Value = [double]"0.001790599999999642"
